# repull data, push all data, mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dSF (column F) values for the affected rows
$ws.Range("F2").Value = -12
$ws.Range("F6").Value = -7
$ws.Range("F9").Value = -2
$ws.Range("F16").Value = 2
